$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New test-case rows (14 & 15) - values first, in the same left-to-right,
#    top-to-bottom order the original author typed them so shared-string
#    indices line up with the target workbook.
# ---------------------------------------------------------------------------

# Row 14 - OPQA-5733
$ws.Range("A14").Value = "OPQA-5733"
$ws.Range("B14").Value = "Verify that to get all entitlements for the user by passing truid"
$ws.Range("C14").Value = "1PENTITLEMENTS"
$ws.Range("D14").Value = "/entitlements/2bd6b996-150e-4b1e-a5c4-c3789237c89b"
$ws.Range("E14").Value = "GET"
$ws.Range("J14").Value = "status=200|| skus=CMC_CIS_04||skus= CMC_AM||skus= CMC_KG||skus= CMC_GE||skus= CMC_TJ||trial_skus=CMC_CIS_04||trial_skus= CMC_AM||trial_skus= CMC_KG||trial_skus= CMC_GE||trial_skus= CMC_TJ||has_trial_skus=true||X-1P-ENT=CMC"

# Row 15 - OPQA-5734
$ws.Range("B15").Value = "Verify that trai skus status value  if user associate with non trail entitlements by passing truid"
$ws.Range("A15").Value = "OPQA-5734"
$ws.Range("C15").Value = "1PENTITLEMENTS"
$ws.Range("D15").Value = "/entitlements/40541757-5531-4006-bbe7-a46ac7ae2d65"
$ws.Range("E15").Value = "GET"
$ws.Range("J15").Value = "status=200|| has_trial_skus=false||X-1P-ENT=CMC"

# ---------------------------------------------------------------------------
# 2. Row heights
# ---------------------------------------------------------------------------
$ws.Rows.Item(14).RowHeight = 180
$ws.Rows.Item(15).RowHeight = 45

# ---------------------------------------------------------------------------
# 3. Borders - apply the existing thin box border used elsewhere in the
#    sheet to the new rows (row 14 fully, row 15 only columns C:K - A15/B15
#    keep their own pre-existing look and L15 is removed below).
# ---------------------------------------------------------------------------
$ws.Range("A14:L14").Borders.LineStyle = 1
$ws.Range("A14:L14").Borders.Weight = 2
$ws.Range("C15:K15").Borders.LineStyle = 1
$ws.Range("C15:K15").Borders.Weight = 2

# ---------------------------------------------------------------------------
# 4. Fill / wrap tweaks so the new rows visually match the rest of the table
# ---------------------------------------------------------------------------
$ws.Range("F14:G15").Interior.ColorIndex = -4142
$ws.Range("G14:G15").WrapText = $true
$ws.Range("H15").WrapText = $true

# ---------------------------------------------------------------------------
# 5. Drop the trailing, unused L column cells (L13 and L15) so they disappear
#    the same way they did in the authored edit.
# ---------------------------------------------------------------------------
$ws.Range("L13").Clear()
$ws.Range("L15").Clear()

# ---------------------------------------------------------------------------
# 6. View state - scroll / selection like the author left it.
# ---------------------------------------------------------------------------
$ws.Range("L2:L19").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 7

Write-Host "edit applied"
